# Update column G ("K") values on Sheet1, rows 2-23, with recalculated
# strike counts (K) replacing the previous "Strike#" derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 3
    3  = 4
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 3
    9  = 1
    10 = 3
    11 = 4
    12 = 4
    13 = 3
    14 = 2
    15 = 6
    16 = 6
    17 = 4
    18 = 5
    19 = 0
    20 = 4
    21 = 1
    22 = 3
    23 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
